$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The project's non-conformities for rows 6-9 (STATUS column F) have been
# resolved, so their STATUS moves from "En proceso" to "Cerrada".
$ws.Range("F6").Value = "Cerrada"
$ws.Range("F7").Value = "Cerrada"
$ws.Range("F8").Value = "Cerrada"
$ws.Range("F9").Value = "Cerrada"

# Rows 8 and 9 previously carried a slightly different (but visually
# identical) cell format than rows 6/7; normalize their number format so
# they end up sharing the same formatting as the rest of the STATUS column.
$ws.Range("F8:F9").NumberFormat = "GENERAL"

# Update the view: scrolled down so row 4 is the first visible row, with
# the active selection on F14.
$ws.Range("F14").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
